$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# F13: was a numeric 0.18 (18%), becomes the literal text "22.5%"
$origFormat = $ws.Range("F13").NumberFormat
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "22.5%"
$ws.Range("F13").NumberFormat = $origFormat

# F14: progress updated from 50% to 83%
$ws.Range("F14").Value = 0.83

# Weekly tracker updates (rows 30-31, week 10 column K)
$ws.Range("K30").Value = 2
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = 2

$ws.Range("F14").Select()
